$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.144801139831543
$ws.Range("B1").Value = 2.569910287857056
$ws.Range("C1").Value = 6.870341777801514
$ws.Range("D1").Value = 2.081766843795776
$ws.Range("E1").Value = 1.229485273361206
